$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10, pushing the existing rows 10-19 down to 11-20.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the new weekly record (same market/category
# metadata as the surrounding rows, new date and volume).
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = [DateTime]"2022-09-23"
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 100112036
$ws.Cells.Item(10, 7).Value = "Caigua"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 20
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 20000
$ws.Cells.Item(10, 13).Value = 20000
$ws.Cells.Item(10, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 1333
$ws.Cells.Item(10, 17).Value = 15
$ws.Cells.Item(10, 18).Value = "Hortaliza"
